$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# 1. Title paragraph: merge "Use Case " + "–" + " " + "Modify pay info" into one run.
$d.Content.Find.Execute("Use Case – Modify pay info", $true, $false, $false, $false, $false, $true, 1, $false, "Use Case – Modify pay info", 2) | Out-Null

# 2. Date paragraph: split "Date: 2019-05-06" into "Date: 2019-05-0" + "9" runs,
#    and move the _GoBack bookmark here (raw OOXML insert so the bookmark lands
#    exactly at the paragraph-end position, which Bookmarks.Add mishandles).
$p = $d.Paragraphs.Item(5)
$start = $p.Range.Start
$end = $p.Range.End
$r = $d.Range($start, $end - 1)
$r.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Date: 2019-05-0</w:t></w:r><w:r><w:t>9</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 3. Priority paragraph: merge "Priority: " + "Medium" into one run.
$d.Content.Find.Execute("Priority: Medium", $true, $false, $false, $false, $false, $true, 1, $false, "Priority: Medium", 2) | Out-Null

# 4. Description paragraph: merge the five runs into one run (also normalizes
#    "f" + "inance department" -> "finance department").
$descText = "A Finance department employee wishes to change an employee’s payment information. The Finance department employee changes the employee’s payment information. Once the changes are made the finance department employee may click the “save changes” button and the new info will be saved in the system."
$d.Content.Find.Execute($descText, $true, $false, $false, $false, $false, $true, 1, $false, $descText, 2) | Out-Null

# 5. Actors list item: merge "Finance department" + " employee" into one run.
$d.Content.Find.Execute("Finance department employee", $true, $false, $false, $false, $false, $true, 1, $false, "Finance department employee", 2) | Out-Null

# 6. Pre-conditions list item: merge "The f" + "inance department" + " employee must be
#    logged into the system." into one run, then drop the old _GoBack bookmark
#    (it was relocated to the Date paragraph in step 2).
$preText = "The finance department employee must be logged into the system."
$d.Content.Find.Execute($preText, $true, $false, $false, $false, $false, $true, 1, $false, $preText, 2) | Out-Null
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
